# Word COM-interop script implementing:
#   "Update documentation on data collection"
#
# Changes (see commit diff):
#  A. Merge ". " + "Download the file..." into a single run (para 2).
#  B. Append " v2" (italic) after "Travel and Activity Diary" (para 4).
#  C. Expand ". Export the file as a CSV..." into
#     ". Make sure to select the right form. Export the data as a CSV..." (para 4).
#  D. Replace the "Go through the data manually..." sentence with
#     "The days tend to be in PST but the times are in UTC. You'll have to
#     change the days manually so that they are in UTC as well. " and move
#     the _GoBack bookmark to sit right before "Otherwise, step 3..." (para 4).
#  E. The _GoBack bookmark that used to sit at the end of para 5 is
#     implicitly removed because it is re-added (same name) at its new
#     location in para 4 — Word bookmarks are unique by name.

$d = $word.ActiveDocument

# Unicode right single quotation mark used throughout the document.
$rsq = [char]0x2019

# --- A. "1. Collect GPS data" paragraph ----------------------------------
# The two runs ". " and "Download the file..." already share identical
# formatting; replacing the second run's own text with itself causes the
# two adjacent runs to coalesce into one, exactly as in the target.
$p2 = $d.Paragraphs(2).Range
$p2.Find.Execute( `
    "Download the file as a tab-delimited text file to your local machine.", `
    $false, $false, $false, $false, $false, $true, 1, $false, `
    "Download the file as a tab-delimited text file to your local machine.", 2) | Out-Null

# --- B. "2. Record Ground Truth" paragraph -------------------------------
# Add " v2" right after "Travel and Activity Diary" (keeps the italic run).
$p4 = $d.Paragraphs(4).Range
$p4.Find.Execute("Activity Diary", $false, $false, $false, $false, $false, `
    $true, 1, $false, "Activity Diary v2", 2) | Out-Null

# --- C. Expand the "Export the file..." sentence -------------------------
$p4 = $d.Paragraphs(4).Range
$p4.Find.Execute("Export the file as a CSV", $false, $false, $false, $false, $false, `
    $true, 1, $false, `
    "Make sure to select the right form. Export the data as a CSV", 2) | Out-Null

# --- D. Replace "Go through the data manually..." sentence ---------------
# Include the single leading (bold) space from the previous run in the
# match so the match does not start exactly on the run boundary (which
# would otherwise mis-attribute formatting); the extra leading space is
# fixed back up to non-bold immediately afterward.
$p4 = $d.Paragraphs(4).Range
$oldSentence = " Go through the data manually to make sure that records " + `
    "corresponding to your data are correct. If they aren" + $rsq + `
    "t, change them so that they are. "
$newSentence = " The days tend to be in PST but the times are in UTC. " + `
    "You" + $rsq + "ll have to change the days manually so that they are " + `
    "in UTC as well. "
$p4.Find.Execute($oldSentence, $false, $false, $false, $false, $false, `
    $true, 1, $false, $newSentence, 2) | Out-Null

# $p4 now spans exactly the replacement text (leading space included);
# restore correct (non-bold) character formatting on everything except
# that first inherited space.
$fixRng = $d.Range($p4.Start + 1, $p4.End)
$fixRng.Font.Bold = $false

# Move the _GoBack bookmark so it sits right before "Otherwise, step 3...".
# Adding a bookmark with a name that is already in use relocates it, so
# the old occurrence at the end of paragraph 5 disappears automatically.
$bmRng = $d.Paragraphs(4).Range
$bmRng.Find.Execute("Otherwise, step 3", $false, $false, $false, $false, $false, `
    $true, 1, $false, "", 0) | Out-Null
$bmRng.Collapse(1)
$d.Bookmarks.Add("_GoBack", $bmRng)

Write-Output "edit complete"
